# Update column F ("dSF") values for specific rows to reflect the
# repulled/recalculated data (repull data, push all data, mean calculation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    15 = 4
    25 = -2
    34 = 4
    36 = -5
    37 = 0
    38 = -2
    46 = -6
    48 = -2
    50 = -2
    51 = -7
    52 = -3
    57 = 0
    59 = -2
    61 = 0
    63 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
